$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (D) and "Volume(1h)" (E) columns with the latest
# crypto snapshot values. Price cells whose text looks like a plain
# number (e.g. "566.36") are written with a leading apostrophe so Excel
# keeps them as text instead of auto-converting them to a numeric value,
# matching the original inline-string cell content.

$ws.Range('D2').Value = '64.854.00'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '2.936.53'
$ws.Range('E3').Value = '  -1.51%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '''566.36'
$ws.Range('E5').Value = '  -2.67%  '
$ws.Range('D6').Value = '''157.11'
$ws.Range('E6').Value = '  +2.25%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '''0.517'
$ws.Range('E8').Value = '  +0.53%  '
$ws.Range('D9').Value = '2.937.97'
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('D10').Value = '''6.67'
$ws.Range('E10').Value = '  -4.69%  '
$ws.Range('D11').Value = '''0.150'
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').Value = '''0.454'
$ws.Range('E12').Value = '  +1.53%  '
$ws.Range('D13').Value = '''0.0000243'
$ws.Range('E13').Value = '  +2.33%  '
$ws.Range('D14').Value = '''33.82'
$ws.Range('E14').Value = '  -0.39%  '
$ws.Range('E15').Value = '  -0.55%  '
$ws.Range('D16').Value = '65.142.32'
$ws.Range('E16').Value = '  +0.32%  '
$ws.Range('D17').Value = '3.430.28'
$ws.Range('E17').Value = '  -1.21%  '
$ws.Range('D18').Value = '''6.91'
$ws.Range('E18').Value = '  +0.07%  '
$ws.Range('D19').Value = '2.943.54'
$ws.Range('E19').Value = '  -1.34%  '
$ws.Range('D20').Value = '''445.14'
$ws.Range('E20').Value = '  -0.77%  '
$ws.Range('D21').Value = '''13.75'
$ws.Range('E21').Value = '  +0.65%  '
$ws.Range('D22').Value = '''0.676'
$ws.Range('E22').Value = '  -0.65%  '
$ws.Range('D23').Value = '''7.18'
$ws.Range('E23').Value = '  -1.48%  '
$ws.Range('D24').Value = '''82.57'
$ws.Range('E24').Value = '  +1.95%  '
$ws.Range('D25').Value = '''2.17'
$ws.Range('E25').Value = '  -1.65%  '
$ws.Range('D26').Value = '''11.96'
$ws.Range('E26').Value = '  -3.43%  '
$ws.Range('D27').Value = '''1.00'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = '''9.95'
$ws.Range('E28').Value = '  -7.36%  '
$ws.Range('D29').Value = '''7.88'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = '''2.33'
$ws.Range('E30').Value = '  -3.21%  '
$ws.Range('D31').Value = '''2.56'
$ws.Range('E31').Value = '  -0.95%  '
$ws.Range('D32').Value = '0.0₃0980'
$ws.Range('E32').Value = '  -4.00%  '
$ws.Range('D33').Value = '''27.11'
$ws.Range('E33').Value = '  +1.29%  '
$ws.Range('E34').Value = '  -1.25%  '
$ws.Range('E35').Value = '  +0.26%  '
$ws.Range('D36').Value = '''0.973'
$ws.Range('E36').Value = '  -0.89%  '
$ws.Range('D37').Value = '''5.64'
$ws.Range('E37').Value = '  -1.22%  '
$ws.Range('D39').Value = '''1.96'
$ws.Range('E39').Value = '  -6.05%  '
$ws.Range('D40').Value = '''43.13'
$ws.Range('E40').Value = '  -4.14%  '
$ws.Range('D41').Value = '''0.295'
$ws.Range('E41').Value = '  -1.65%  '
$ws.Range('D42').Value = '''0.118'
$ws.Range('E42').Value = '  -1.39%  '
$ws.Range('D43').Value = '''2.77'
$ws.Range('E43').Value = '  -4.49%  '
$ws.Range('D44').Value = '''8.38'
$ws.Range('E44').Value = '  -0.03%  '
$ws.Range('D45').Value = '''379.68'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('D46').Value = '''0.0349'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').Value = '2.722.14'
$ws.Range('E47').Value = '  -1.68%  '
$ws.Range('D48').Value = '''131.75'
$ws.Range('E48').Value = '  -1.65%  '
$ws.Range('D50').Value = '''0.106'
$ws.Range('E50').Value = '  +0.78%  '
$ws.Range('D51').Value = '''2.13'
$ws.Range('E51').Value = '  +4.34%  '
